# Weekly update: insert two new "Locoto" price records into the time
# series, shifting the existing rows down (newest-first ordering).
#
#   - a new record is inserted at row 10 (date 2023-02-13 / serial 44970)
#   - a new record is inserted at row 17 (date 2023-02-14 / serial 44971)
#
# Everything else (rows 2-9, and the data that used to live at rows
# 10-65) shifts down accordingly; the sheet grows from 65 to 67 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new record at row 10 --------------------------------
$ws.Rows("10").Insert()

$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 44970
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 100112042
$ws.Range("G10").Value = "Locoto"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 140
$ws.Range("K10").Value = 3300
$ws.Range("L10").Value = 3300
$ws.Range("M10").Value = 3300
$ws.Range("N10").Value = "`$/kilo"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 3300
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"

# --- Insert the second new record at row 17 --------------------------------
$ws.Rows("17").Insert()

$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "Vega Modelo de Temuco"
$ws.Range("C17").Value = "La Araucanía"
$ws.Range("D17").Value = 44971
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = 100112042
$ws.Range("G17").Value = "Locoto"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 3300
$ws.Range("L17").Value = 3300
$ws.Range("M17").Value = 3300
$ws.Range("N17").Value = "`$/kilo"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 3300
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = "Hortaliza"
